# Add a "Team Cool Kids" heading paragraph right after the title
# ("Business Requirements Document Template") and move the document's
# "_GoBack" bookmark (last-edit marker) onto the new text, exactly as
# Word itself would after a user typed the new line.

$d = $word.ActiveDocument

# --- 1. Insert a new paragraph right after the title paragraph ---------
$titlePara = $d.Paragraphs.Item(1)
$titleRange = $titlePara.Range
$titleRange.InsertParagraphAfter()

$newPara = $d.Paragraphs.Item(2)
$newRange = $newPara.Range

# Type the team name, plus a throw-away trailing marker character that
# we use purely to anchor the bookmark precisely; it is removed again
# right after. (Placing a zero-length bookmark exactly at the end of a
# paragraph's text can't be done directly, so we park it before the
# marker character and then delete that character.)
$newRange.InsertAfter("Team Cool KidsZ")

$newPara2 = $d.Paragraphs.Item(2)
$newRange2 = $newPara2.Range
$markerPos = $newRange2.End - 2
$bookmarkSpot = $d.Range($markerPos, $markerPos)

# --- 2. Move the "_GoBack" bookmark from wherever it currently lives ---
#        to sit right after "Team Cool Kids" in the new paragraph -------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}
$d.Bookmarks.Add("_GoBack", $bookmarkSpot)

# --- 3. Remove the temporary marker character ---------------------------
$markerRange = $d.Range($markerPos, $markerPos + 1)
$markerRange.Delete()
